# Update data/markov to newest formats
# Adds a small "Format: v0.1.0" tag in row 2 of the Global Parameters sheet,
# using a new italic font style (right-aligned label, left value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global Parameters")

# Label cell B2: "Format:"
$b2 = $ws.Range("B2")
$b2.Value = "Format:"
$b2.Font.Italic = $true
$b2.Font.Name = "Aptos"
$b2.Font.Size = 11
$b2.HorizontalAlignment = -4152   # xlRight
$b2.VerticalAlignment = -4108     # xlCenter

# Value cell C2: "v0.1.0"
$c2 = $ws.Range("C2")
$c2.Value = "v0.1.0"
$c2.Font.Italic = $true
$c2.Font.Name = "Aptos"
$c2.Font.Size = 11
$c2.VerticalAlignment = -4108     # xlCenter

$ws.Rows.Item(2).RowHeight = 18.75
